$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Linde-Decatur'
$ws.Range("B2").Value = 'Decatur'
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 39
$ws.Range("E2").Value = 196.63113068651
$ws.Range("F2").Value = 56549
$ws.Range("G2").Value = 287.589253047405
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.1
$ws.Range("J2").Value = 37934
$ws.Range("K2").Value = 13220
$ws.Range("L2").Value = 293
$ws.Range("M2").Value = 165
$ws.Range("N2").Value = 7775
$ws.Range("O2").Value = 49.0047692307692
$ws.Range("P2").Value = 6.9644053965162
$ws.Range("Q2").Value = 12.9274906797582
$ws.Range("R2").Value = 40
$ws.Range("S2").Value = 0.492307692307692

$ws.Range("A3").Value = 'A-GAS'
$ws.Range("B3").Value = 'El Dorado'
$ws.Range("C3").Value = 106401
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 422.192389752563
$ws.Range("F3").Value = 8797
$ws.Range("G3").Value = 20.8364722186388
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 7211
$ws.Range("K3").Value = 1335
$ws.Range("L3").Value = 173
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 298
$ws.Range("O3").Value = 55.6223333333333
$ws.Range("P3").Value = 4.00529428733787
$ws.Range("Q3").Value = 11.1294509931195
$ws.Range("R3").Value = 50
$ws.Range("S3").Value = 0.566666666666667

$ws.Range("A4").Value = 'AEROPRES-SANDIMAS'
$ws.Range("B4").Value = 'San Dimas'
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 206
$ws.Range("E4").Value = 410.139393764204
$ws.Range("F4").Value = 365058
$ws.Range("G4").Value = 890.082751255731
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.039
$ws.Range("J4").Value = 210194
$ws.Range("K4").Value = 14212
$ws.Range("L4").Value = 4746
$ws.Range("M4").Value = 52470
$ws.Range("N4").Value = 181109
$ws.Range("O4").Value = 85.066345177665
$ws.Range("P4").Value = 5.21839427599669
$ws.Range("Q4").Value = 6.06202061612871
$ws.Range("R4").Value = 37.9611650485437
$ws.Range("S4").Value = 0.440291262135922

$ws.Range("A5").Value = 'CALAMCO'
$ws.Range("B5").Value = 'Stockton'
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 164
$ws.Range("E5").Value = 279.208811816867
$ws.Range("F5").Value = 273798
$ws.Range("G5").Value = 980.620913137887
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.024
$ws.Range("J5").Value = 141614
$ws.Range("K5").Value = 27639
$ws.Range("L5").Value = 1867
$ws.Range("M5").Value = 43305
$ws.Range("N5").Value = 137440
$ws.Range("O5").Value = 49.7378944099379
$ws.Range("P5").Value = 9.0299994623911
$ws.Range("Q5").Value = 12.0943532838632
$ws.Range("R5").Value = 30
$ws.Range("S5").Value = 0.479878048780488

$ws.Range("A6").Value = 'Linde-Whiting'
$ws.Range("B6").Value = 'East Chicago'
$ws.Range("C6").Value = 183562
$ws.Range("D6").Value = 165
$ws.Range("E6").Value = 80.1469498240182
$ws.Range("F6").Value = 178959
$ws.Range("G6").Value = 2232.88597249112
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 82166
$ws.Range("K6").Value = 57189
$ws.Range("L6").Value = 672
$ws.Range("M6").Value = 1248
$ws.Range("N6").Value = 68444
$ws.Range("O6").Value = 45.8696233766234
$ws.Range("P6").Value = 10.6961310497243
$ws.Range("Q6").Value = 11.7969216170977
$ws.Range("R6").Value = 30
$ws.Range("S6").Value = 0.364242424242424

$ws.Range("A7").Value = 'Diversified-CPC'
$ws.Range("B7").Value = 'Channahon'
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 248.778205088469
$ws.Range("F7").Value = 54788
$ws.Range("G7").Value = 220.228295242007
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.52
$ws.Range("J7").Value = 47218
$ws.Range("K7").Value = 4051
$ws.Range("L7").Value = 53
$ws.Range("M7").Value = 1065
$ws.Range("N7").Value = 8657
$ws.Range("O7").Value = 86.4332
$ws.Range("P7").Value = 3.30740926855136
$ws.Range("Q7").Value = 3.26505024426427
$ws.Range("R7").Value = 27
$ws.Range("S7").Value = 0.335

$ws.Range("A8").Value = 'CFI-PortNeal'
$ws.Range("B8").Value = 'Sergeant Bluff'
$ws.Range("C8").Value = 2888938
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 388.310293949375
$ws.Range("F8").Value = 14583
$ws.Range("G8").Value = 37.5550177969303
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0.9
$ws.Range("J8").Value = 11581
$ws.Range("K8").Value = 37
$ws.Range("L8").Value = 2151
$ws.Range("M8").Value = 203
$ws.Range("N8").Value = 1008
$ws.Range("O8").Value = 74.213
$ws.Range("P8").Value = 4.57628419421131
$ws.Range("Q8").Value = 7.46600688747593
$ws.Range("R8").Value = 20
$ws.Range("S8").Value = 0.2

$ws.Range("A9").Value = 'APC-Geismar'
$ws.Range("B9").Value = 'Geismar'
$ws.Range("C9").Value = 56325
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = 180.847413422864
$ws.Range("F9").Value = 44136
$ws.Range("G9").Value = 244.051043720485
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0.71
$ws.Range("J9").Value = 24770
$ws.Range("K9").Value = 17026
$ws.Range("L9").Value = 18
$ws.Range("M9").Value = 436
$ws.Range("N9").Value = 2072
$ws.Range("O9").Value = 63.8595882352941
$ws.Range("P9").Value = 10.0047370993105
$ws.Range("Q9").Value = 7.22680275103533
$ws.Range("R9").Value = 88.8235294117647
$ws.Range("S9").Value = 0.517647058823529

$ws.Range("A10").Value = 'Honeywell-Geismar'
$ws.Range("B10").Value = 'Geismar'
$ws.Range("C10").Value = 528957
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 228.365522556404
$ws.Range("F10").Value = 38167
$ws.Range("G10").Value = 167.131183257197
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0.93
$ws.Range("J10").Value = 23475
$ws.Range("K10").Value = 13551
$ws.Range("L10").Value = 26
$ws.Range("M10").Value = 199
$ws.Range("N10").Value = 1091
$ws.Range("O10").Value = 67.4632307692308
$ws.Range("P10").Value = 6.08715442194185
$ws.Range("Q10").Value = 3.24084657418216
$ws.Range("R10").Value = 115.714285714286
$ws.Range("S10").Value = 0.55

$ws.Range("A11").Value = 'APC-PortAuthur'
$ws.Range("B11").Value = 'Port Arthur'
$ws.Range("C11").Value = 2423361
$ws.Range("D11").Value = 51
$ws.Range("E11").Value = 956.159882202067
$ws.Range("F11").Value = 52463
$ws.Range("G11").Value = 54.8684388213152
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.14
$ws.Range("J11").Value = 26995
$ws.Range("K11").Value = 19176
$ws.Range("L11").Value = 257
$ws.Range("M11").Value = 3777
$ws.Range("N11").Value = 18545
$ws.Range("O11").Value = 39.1885434782609
$ws.Range("P11").Value = 12.1854389916085
$ws.Range("Q11").Value = 12.5620748615726
$ws.Range("R11").Value = 46.8627450980392
$ws.Range("S11").Value = 0.398039215686275

$ws.Range("A12").Value = 'AEROPRES-SIBLEY'
$ws.Range("B12").Value = 'Sibley'
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = 21
$ws.Range("E12").Value = 234.617951559258
$ws.Range("F12").Value = 21636
$ws.Range("G12").Value = 92.218007429561
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0.52
$ws.Range("J12").Value = 12191
$ws.Range("K12").Value = 8911
$ws.Range("L12").Value = 52
$ws.Range("M12").Value = 104
$ws.Range("N12").Value = 243
$ws.Range("O12").Value = 30.7228095238095
$ws.Range("P12").Value = 8.50433306877584
$ws.Range("Q12").Value = 23.8404653935191
$ws.Range("R12").Value = 40
$ws.Range("S12").Value = 0.5

$ws.Range("A13").Value = 'HaltermanCarless'
$ws.Range("B13").Value = 'Manvel'
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = 28
$ws.Range("E13").Value = 178.547661387283
$ws.Range("F13").Value = 81631
$ws.Range("G13").Value = 457.194450858341
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0.5
$ws.Range("J13").Value = 57444
$ws.Range("K13").Value = 14188
$ws.Range("L13").Value = 387
$ws.Range("M13").Value = 4159
$ws.Range("N13").Value = 27398
$ws.Range("O13").Value = 73.4291153846154
$ws.Range("P13").Value = 5.56472923238347
$ws.Range("Q13").Value = 6.55266548618066
$ws.Range("R13").Value = 30
$ws.Range("S13").Value = 0.367857142857143

$ws.Range("A14").Value = 'Chemours-CorpusChristie'
$ws.Range("B14").Value = 'Gregory'
$ws.Range("C14").Value = 26401
$ws.Range("D14").Value = 21
$ws.Range("E14").Value = 334.639503004682
$ws.Range("F14").Value = 34683
$ws.Range("G14").Value = 103.642874462178
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.76
$ws.Range("J14").Value = 31808
$ws.Range("K14").Value = 608
$ws.Range("L14").Value = 183
$ws.Range("M14").Value = 485
$ws.Range("N14").Value = 15487
$ws.Range("O14").Value = 69.3915789473684
$ws.Range("P14").Value = 4.32173615075584
$ws.Range("Q14").Value = 2.65422327250218
$ws.Range("R14").Value = 20
$ws.Range("S14").Value = 0.20952380952381

$ws.Range("A15").Value = 'DiversifiedG&O'
$ws.Range("B15").Value = 'Tad'
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 195.038768642775
$ws.Range("F15").Value = 12136
$ws.Range("G15").Value = 62.2235265555218
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 11672
$ws.Range("K15").Value = 117
$ws.Range("L15").Value = 6
$ws.Range("M15").Value = 195
$ws.Range("N15").Value = 109
$ws.Range("O15").Value = 48.6558888888889
$ws.Range("P15").Value = 6.55367935362492
$ws.Range("Q15").Value = 11.4346719252667
$ws.Range("R15").Value = 30
$ws.Range("S15").Value = 0.4
